# Roster update: add "Keyonte George" at the top of the list (row 2) and
# re-sequence the remaining players, keeping their Position/Team pairs intact.
# The final row order (Name, Position, Team) is:

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Keyonte George",        "PG,SG",    "Utah Jazz"),
    @("Anthony Edwards",       "SG,SF",    "Minnesota Timberwolves"),
    @("James Harden",          "PG,SG",    "LA Clippers"),
    @("Anfernee Simons",       "PG,SG",    "Portland Trail Blazers"),
    @("Jayson Tatum",          "SF,PF",    "Boston Celtics"),
    @("Paul George",           "SG,SF,PF", "Philadelphia 76ers"),
    @("Amen Thompson",         "SG,SF,PF", "Houston Rockets"),
    @("Zion Williamson",       "PF,C",     "New Orleans Pelicans"),
    @("Nicolas Claxton",       "C",        "Brooklyn Nets"),
    @("Zach Edey",             "C",        "Memphis Grizzlies"),
    @("Jaren Jackson Jr.",     "PF,C",     "Memphis Grizzlies"),
    @("Ivica Zubac",           "C",        "LA Clippers"),
    @("Bobby Portis",          "PF,C",     "Milwaukee Bucks"),
    @("Draymond Green",        "PF,C",     "Golden State Warriors"),
    @("RJ Barrett",            "SG,SF,PF", "Toronto Raptors"),
    @("Giannis Antetokounmpo", "PF,C",     "Milwaukee Bucks"),
    @("Fred VanVleet",         "PG",       "Houston Rockets"),
    @("Bradley Beal",          "PG,SG,SF", "Phoenix Suns")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
